$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 170098
$ws.Range("C4").Value = 160926
$ws.Range("C7").Value = 5.39
$ws.Range("C8").Value = 65.73999999999999
